# Data-driven application of the scheduled market-data refresh.
# Each entry: Sheet name, Row number, then values for columns H..N
# (use $null to mean 'clear this cell entirely').
$updates = @(
    [PSCustomObject]@{ Sheet = "ALC"; Row = 51; Vals = @(3999.6667, 4000, 3999.5, 4000, 3999.5, -3516, -4967.5) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 62; Vals = @(16764.072, 14782.833, 18250, 14782.833, 18250, -14158.833, -19498) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 63; Vals = @(80271, 0, 80271, 0, 80271, $null, -81519) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 65; Vals = @(16764.072, 14782.833, 18250, 73914.16500000001, 91250, -70794.16500000001, -97490) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 66; Vals = @(80271, 0, 80271, 0, 240813, $null, -247053) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 76; Vals = @(4522.6313, 3319.25, 6585.5713, 3319.25, 6585.5713, -3004.25, -7215.5713) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 79; Vals = @(4522.6313, 3319.25, 6585.5713, 3319.25, 6585.5713, -2227.25, -8769.5713) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 86; Vals = @(4565.4287, 2166.6667, 5219.636, 2166.6667, 5219.636, -1043.6667, -7465.636) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 89; Vals = @(4565.4287, 2166.6667, 5219.636, 10833.3335, 26098.18, -5217.333500000001, -37330.18) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 98; Vals = @(2802.111, 2866.647, 1705, 2866.647, 1705, -1368.647, -4701) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 106; Vals = @(2245, 1490, 3000, 1490, 3000, -859, -4262) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 113; Vals = @(1445.238, 1489.9231, 1372.625, 1489.9231, 1372.625, 1764.0769, -7880.625) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 122; Vals = @(2802.111, 2866.647, 1705, 8599.940999999999, 5115, -6149.940999999999, -10015) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 127; Vals = @(4175, 4824, 930, 14472, 2790, -9512, -12710) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 137; Vals = @(1049591.5, 1417, 1821930.6, 4251, 5465791.800000001, -1701, -5470891.800000001) }
    [PSCustomObject]@{ Sheet = "ALC"; Row = 138; Vals = @(3222.6667, 2788.9092, 3473.7896, 8366.7276, 10421.3688, -3226.7276, -20701.3688) }
    [PSCustomObject]@{ Sheet = "ARM"; Row = 32; Vals = @(5214644.5, 5469037, 41994.332, 5469037, 41994.332, -5468750, -42568.332) }
    [PSCustomObject]@{ Sheet = "ARM"; Row = 122; Vals = @(3764.6843, 2289.8572, 4625, 6869.571599999999, 13875, -4419.571599999999, -18775) }
    [PSCustomObject]@{ Sheet = "BSM"; Row = 20; Vals = @(1899, 1828.2941, 2500, 1828.2941, 2500, -1581.2941, -2994) }
    [PSCustomObject]@{ Sheet = "BSM"; Row = 22; Vals = @(1812.25, 1928.2858, 1000, 1928.2858, 1000, -1755.2858, -1346) }
    [PSCustomObject]@{ Sheet = "CRP"; Row = 25; Vals = @(649.75, 649.75, 0, 649.75, 0, -475.75, $null) }
    [PSCustomObject]@{ Sheet = "CRP"; Row = 31; Vals = @(6666.079, 1751.4286, 9532.958000000001, 1751.4286, 9532.958000000001, -1456.4286, -10122.958) }
    [PSCustomObject]@{ Sheet = "CRP"; Row = 34; Vals = @(6666.079, 1751.4286, 9532.958000000001, 1751.4286, 9532.958000000001, -1549.4286, -9936.958000000001) }
    [PSCustomObject]@{ Sheet = "CRP"; Row = 41; Vals = @(37249.082, 19750, 45998.625, 19750, 45998.625, -19322, -46854.625) }
    [PSCustomObject]@{ Sheet = "CRP"; Row = 50; Vals = @(38899.6, 8749.5, 58999.668, 8749.5, 58999.668, -8124.5, -60249.668) }
    [PSCustomObject]@{ Sheet = "CRP"; Row = 51; Vals = @(43749.75, 25000, 49999.668, 25000, 49999.668, -24264, -51471.668) }
    [PSCustomObject]@{ Sheet = "CRP"; Row = 59; Vals = @(53484.25, 104, 71277.664, 104, 71277.664, 1041, -73567.664) }
    [PSCustomObject]@{ Sheet = "CRP"; Row = 60; Vals = @(22839.8, 19666.334, 24199.857, 19666.334, 24199.857, -19155.334, -25221.857) }
    [PSCustomObject]@{ Sheet = "CRP"; Row = 61; Vals = @(43749.75, 25000, 49999.668, 25000, 49999.668, -24652, -50695.668) }
    [PSCustomObject]@{ Sheet = "CRP"; Row = 107; Vals = @(552.6316, 558.82355, 500, 558.82355, 500, 1361.17645, -4340) }
    [PSCustomObject]@{ Sheet = "CUL"; Row = 4; Vals = @(39738720, 45056908, 6500037, 135170724, 19500111, -135170612, -19500335) }
    [PSCustomObject]@{ Sheet = "CUL"; Row = 11; Vals = @(12500840, 371, 14286622, 1113, 42859866, -973, -42860146) }
    [PSCustomObject]@{ Sheet = "CUL"; Row = 132; Vals = @(2859.9524, 2154.5715, 3212.6428, 19391.1435, 28913.7852, -16861.1435, -33973.7852) }
    [PSCustomObject]@{ Sheet = "CUL"; Row = 141; Vals = @(0, 0, 0, 0, 0, $null, $null) }
    [PSCustomObject]@{ Sheet = "GSM"; Row = 10; Vals = @(5000, 0, 5000, 0, 5000, $null, -5338) }
    [PSCustomObject]@{ Sheet = "GSM"; Row = 11; Vals = @(80000, 80000, 0, 80000, 0, -79861, $null) }
    [PSCustomObject]@{ Sheet = "GSM"; Row = 21; Vals = @(690000, 10000000, 25000, 10000000, 25000, -9999827, -25346) }
    [PSCustomObject]@{ Sheet = "GSM"; Row = 30; Vals = @(690000, 10000000, 25000, 10000000, 25000, -9999895, -25210) }
    [PSCustomObject]@{ Sheet = "GSM"; Row = 44; Vals = @(40000, 0, 40000, 0, 40000, $null, -41192) }
    [PSCustomObject]@{ Sheet = "GSM"; Row = 105; Vals = @(90180.664, 0, 90180.664, 0, 90180.664, $null, -97168.664) }
    [PSCustomObject]@{ Sheet = "GSM"; Row = 113; Vals = @(2163.0967, 1312.0667, 2960.9375, 1312.0667, 2960.9375, 857.9332999999999, -7300.9375) }
    [PSCustomObject]@{ Sheet = "GSM"; Row = 135; Vals = @(64994.8, 0, 64994.8, 0, 64994.8, $null, -75134.8) }
    [PSCustomObject]@{ Sheet = "LTW"; Row = 7; Vals = @(6737.1113, 7206.8, 6150, 7206.8, 6150, -7094.8, -6374) }
    [PSCustomObject]@{ Sheet = "LTW"; Row = 22; Vals = @(54682.367, 250737.25, 2401.0667, 250737.25, 2401.0667, -250442.25, -2991.0667) }
    [PSCustomObject]@{ Sheet = "LTW"; Row = 23; Vals = @(1600, 1600, 0, 1600, 0, -1370, $null) }
    [PSCustomObject]@{ Sheet = "LTW"; Row = 27; Vals = @(54682.367, 250737.25, 2401.0667, 250737.25, 2401.0667, -250630.25, -2615.0667) }
    [PSCustomObject]@{ Sheet = "LTW"; Row = 126; Vals = @(6737.1113, 7206.8, 6150, 21620.4, 18450, -19150.4, -23390) }
    [PSCustomObject]@{ Sheet = "WVR"; Row = 126; Vals = @(3198.1538, 2658.8572, 3827.3333, 7976.571599999999, 11481.9999, -5506.571599999999, -16421.9999) }
)

$wb = $excel.ActiveWorkbook
$cols = @("H", "I", "J", "K", "L", "M", "N")

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $u.Vals[$i]
        $cell = $ws.Range("$col$($u.Row)")
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}

Write-Host "Applied $($updates.Count) row updates across $((($updates | Select-Object -ExpandProperty Sheet -Unique)).Count) sheets."
